# Feature: add arrows (arrow_n).
#
# The "meta" sheet stores key/value pairs (one pair per row, key in column A,
# value in column B) and is terminated by a trailing, otherwise-empty row
# that only carries column A's style. This adds a new "style" / "default"
# metadata pair right before that trailing placeholder row, pushing the
# placeholder row down by one.

$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("meta")

# The last populated key/value row before the trailing placeholder row.
$lastKeyRow = 5
$placeholderRow = $lastKeyRow + 1

# Insert a new, blank row above the trailing placeholder row so the
# placeholder (and its formatting) moves down one row, exactly as it did
# before this edit.
$meta.Rows.Item($placeholderRow).Insert(-4121) | Out-Null

# Fill in the newly inserted row with the new metadata entry.
$meta.Range("A" + $placeholderRow).Value = "style"
$meta.Range("B" + $placeholderRow).Value = "default"
